# Updates the cryptos list data (columns B/C/D/E) to reflect the latest
# scraped values, matching the GitHub Actions "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch keep their existing "Text" representation so
# that values like "42.437.14", "0.800" or "  +1.47%  " are not reinterpreted
# by Excel as numbers/dates and lose formatting (leading/trailing spaces,
# trailing zeros, multiple dot separators, etc.).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.437.14'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = '2.288.63'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '156.94'
$ws.Range("E5").Value = '  +15,583.12%  '
$ws.Range("D6").Value = '307.13'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '95.78'
$ws.Range("E7").Value = '  +4.54%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '0.495'
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").Value = '35.99'
$ws.Range("E11").Value = '  +11.10%  '
$ws.Range("D12").Value = '0.0804'
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").Value = '2.643.26'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").Value = '2.295.73'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '0.800'
$ws.Range("E18").Value = '  +5.11%  '
$ws.Range("D19").Value = '42.310.19'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '12.68'
$ws.Range("E20").Value = '  +3.82%  '
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("E22").Value = '  +1.63%  '
$ws.Range("D23").Value = '68.14'
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").Value = '243.17'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = '24.13'
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = '36.26'
$ws.Range("E29").Value = '  +5.95%  '
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").Value = '2.10'
$ws.Range("E31").Value = '  -8.86%  '
$ws.Range("D32").Value = '161.86'
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("E36").Value = '  +2.97%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.108'
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '17.25'
$ws.Range("E38").Value = '  +2.80%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  +7.62%  '
$ws.Range("D43").Value = '2.014.09'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").Value = '19.54'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("E45").Value = '  +10.81%  '
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("D47").Value = '10.17'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '2.99'
$ws.Range("E48").Value = '  +4.15%  '
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("D51").Value = '73.09'
$ws.Range("E51").Value = '  -0.64%  '
